$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New cell content (order matters: it controls the order new entries are
#    appended to the shared-strings table, which must match the target file).
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "theoretical"
$ws.Range("B6").Value = "/θɪəˈret.ɪ.kəl/"
$ws.Range("D6").Value = "理论上的"
$ws.Range("C1").Value = "part of speech"
$ws.Range("D1").Value = "meaning"
$ws.Range("C2").Value = "[n.] [v.]"

# ---------------------------------------------------------------------------
# 2. Column widths (column B grew wider, new column C introduced).
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 16.7
$ws.Columns("C").ColumnWidth = 15

# ---------------------------------------------------------------------------
# 3. Borders - apply a thin box border to every cell in A1:D30. Clear any
#    pre-existing border first so a stray leftover (B12) cannot pollute the
#    new uniform border style, then apply per full-height column so every
#    cell individually resolves to "thin" on all four sides.
# ---------------------------------------------------------------------------
$ws.Range("A1:D30").Borders.LineStyle = -4142

$ws.Range("A1:A30").Borders.LineStyle = 1
$ws.Range("A1:A30").Borders.Weight = 2

$ws.Range("B1:B30").Borders.LineStyle = 1
$ws.Range("B1:B30").Borders.Weight = 2

$ws.Range("C1:C30").Borders.LineStyle = 1
$ws.Range("C1:C30").Borders.Weight = 2

$ws.Range("D1:D30").Borders.LineStyle = 1
$ws.Range("D1:D30").Borders.Weight = 2

# ---------------------------------------------------------------------------
# 4. Page setup - switched to portrait / paper size 9, page-break preview.
# ---------------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

$win = $excel.ActiveWindow
$win.View = 3

# ---------------------------------------------------------------------------
# 5. Selection.
# ---------------------------------------------------------------------------
$ws.Range("C5").Select() | Out-Null
